$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rates text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 11.46 = 46201.1 pesos`n✅ 46201.1 pesos = 11.41 = 966.23 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 87.249
$wsTasas.Range("O10").Value = 4031
$wsTasas.Range("N12").Value = 4050
$wsTasas.Range("O12").Value = 84.7
